# Generate Report for Handoff
# Update status + timestamps for the localization-status report and widen
# the "Status" column on each sheet to fit the new "Ready for handoff" text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Target column width is ~17.22 characters (to fit "Ready for handoff").
# ColumnWidth is quantized by the host to 1/6-character steps, so feed it a
# value from the middle of the input bucket that rounds to the width closest
# to the desired 17.2159881591797.
$statusColWidth = 16.333333333333336

# --- Overview sheet ---------------------------------------------------
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-09-03 12:43:09"

$overview.Columns.Item(5).ColumnWidth = $statusColWidth
$overview.Columns.Item(6).ColumnWidth = $statusColWidth

# --- zh-cn sheet --------------------------------------------------------
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-09-03 12:43:00"

$zhcn.Columns.Item(3).ColumnWidth = $statusColWidth

# --- de-de sheet --------------------------------------------------------
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-09-03 12:43:09"

$dede.Columns.Item(3).ColumnWidth = $statusColWidth
